$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text so numeric-looking strings (e.g. "1.002") are
# stored as literal text, matching the original inlineStr cells, instead
# of being auto-coerced to numbers by Excel.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

# --- Price (column D) updates ---
$ws.Range("D2").Value = "24.876.32"
$ws.Range("D3").Value = "1.708.64"
$ws.Range("D4").Value = "1.002"
$ws.Range("D5").Value = "310.43"
$ws.Range("D6").Value = "0.9965"
$ws.Range("D7").Value = "0.3734"
$ws.Range("D8").Value = "49.44"
$ws.Range("D9").Value = "0.3437"
$ws.Range("D10").Value = "1.204"
$ws.Range("D11").Value = "0.07517"
$ws.Range("D12").Value = "0.9984"
$ws.Range("D13").Value = "21.11"
$ws.Range("D14").Value = "6.303"
$ws.Range("D15").Value = "7.014"
$ws.Range("D16").Value = "1.711.58"
$ws.Range("D17").Value = "0.00001133"
$ws.Range("D18").Value = "0.06748"
$ws.Range("D19").Value = "0.9968"
$ws.Range("D20").Value = "84.55"
$ws.Range("D21").Value = "17.28"
$ws.Range("D22").Value = "6.384"
$ws.Range("D23").Value = "13.14"
$ws.Range("D24").Value = "24.855.79"
$ws.Range("D25").Value = "2.432"
$ws.Range("D26").Value = "2.783"
$ws.Range("D27").Value = "20.37"
$ws.Range("D28").Value = "150.86"
$ws.Range("D29").Value = "132.21"
$ws.Range("D30").Value = "1.900.41"
$ws.Range("D31").Value = "1.253"
$ws.Range("D32").Value = "6.943"
$ws.Range("D33").Value = "4.197"
$ws.Range("D34").Value = "1.844"
$ws.Range("D35").Value = "13.84"
$ws.Range("D36").Value = "0.08806"
$ws.Range("D37").Value = "5.598"
$ws.Range("D38").Value = "0.06651"
$ws.Range("D41").Value = "0.2242"
$ws.Range("D42").Value = "1.275"
$ws.Range("D43").Value = "0.6457"
$ws.Range("D45").Value = "13.99"
$ws.Range("D46").Value = "0.6179"
$ws.Range("D47").Value = "3.827"
$ws.Range("D48").Value = "2.139"
$ws.Range("D49").Value = "130.48"
$ws.Range("D50").Value = "0.07322"
$ws.Range("D51").Value = "79.82"

# Row 39/40 swap: VeChain <-> FraxShare traded places with refreshed data.
$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").Value = "9.138"
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "0.02401"

# Restore column D to the workbook default (no explicit number format),
# same as the surrounding untouched cells.
$dRange.Style = "Normal"

# --- Volume(1h) (column E) updates ---
$ws.Range("E2").Value = "  +2.04%  "
$ws.Range("E3").Value = "  +1.74%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("E5").Value = "  +0.88%  "
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("E7").Value = "  +0.86%  "
$ws.Range("E8").Value = "  +3.53%  "
$ws.Range("E9").Value = "  -0.19%  "
$ws.Range("E10").Value = "  +1.85%  "
$ws.Range("E11").Value = "  +3.58%  "
$ws.Range("E12").Value = "  -0.07%  "
$ws.Range("E13").Value = "  +3.40%  "
$ws.Range("E14").Value = "  +3.24%  "
$ws.Range("E15").Value = "  +3.86%  "
$ws.Range("E16").Value = "  +1.68%  "
$ws.Range("E17").Value = "  +1.89%  "
$ws.Range("E18").Value = "  +0.42%  "
$ws.Range("E19").Value = "  +0.02%  "
$ws.Range("E20").Value = "  +4.07%  "
$ws.Range("E21").Value = "  +4.92%  "
$ws.Range("E22").Value = "  +4.59%  "
$ws.Range("E23").Value = "  +9.85%  "
$ws.Range("E24").Value = "  +2.15%  "
$ws.Range("E25").Value = "  -0.14%  "
$ws.Range("E26").Value = "  +4.42%  "
$ws.Range("E27").Value = "  +4.17%  "
$ws.Range("E28").Value = "  -0.63%  "
$ws.Range("E29").Value = "  +3.90%  "
$ws.Range("E30").Value = "  +1.80%  "
$ws.Range("E31").Value = "  +29.07%  "
$ws.Range("E32").Value = "  +10.38%  "
$ws.Range("E33").Value = "  +3.96%  "
$ws.Range("E34").Value = "  +4.87%  "
$ws.Range("E35").Value = "  +12.84%  "
$ws.Range("E36").Value = "  +4.21%  "
$ws.Range("E37").Value = "  +4.69%  "
$ws.Range("E38").Value = "  +3.96%  "
$ws.Range("E39").Value = "  +1.15%  "
$ws.Range("E40").Value = "  +3.12%  "
$ws.Range("E41").Value = "  +6.30%  "
$ws.Range("E42").Value = "  +1.19%  "
$ws.Range("E43").Value = "  +4.73%  "
$ws.Range("E44").Value = "  -0.03%  "
$ws.Range("E45").Value = "  +6.68%  "
$ws.Range("E46").Value = "  +4.14%  "
$ws.Range("E47").Value = "  +1.14%  "
$ws.Range("E48").Value = "  +5.57%  "
$ws.Range("E49").Value = "  +2.41%  "
$ws.Range("E50").Value = "  +1.50%  "
$ws.Range("E51").Value = "  +5.34%  "
